# Teyvat Archives Schema workbook update
# - Adds the "average_rating" / DECIMAL(3,2) field to the Character sheet's
#   Field/Data Type table (pushing the existing created_at / image_url rows
#   down by one row).
# - Normalizes the User sheet's blank E7 cell formatting to match D7.
# - Updates the remembered selections on the Build and Character sheets and
#   makes "Character" the active tab (API Workflow + Schema additions).

$wb = $excel.ActiveWorkbook

# --- User sheet: E7 adopts D7's (border/fill) formatting -------------------
$wsUser = $wb.Worksheets.Item("User")
$wsUser.Range("D7").Copy($wsUser.Range("E7"))

# --- Character sheet: insert the new "average_rating" row ------------------
$wsChar = $wb.Worksheets.Item("Character")

# Shift the existing D8:E9 rows down one row (bottom-up so nothing is
# clobbered before it's been copied), then write the new field in the
# vacated D8:E8 slot.
$wsChar.Range("D9:E9").Copy($wsChar.Range("D10"))
$wsChar.Range("D8:E8").Copy($wsChar.Range("D9"))
$wsChar.Range("D8").Value = "average_rating"
$wsChar.Range("E8").Value = "DECIMAL(3,2)"

# --- Build sheet: remembered selection moves to D6:E6 -----------------------
$wsBuild = $wb.Worksheets.Item("Build")
$wsBuild.Range("D6:E6").Select()

# --- Character becomes the active sheet/tab, selection on E11 --------------
$wsChar.Activate()
$wsChar.Range("E11").Select()
